$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comparisons")

# 1. Header rename for the comparison-algorithm columns
$ws.Range("J2").Value = "ZEUS"
$ws.Range("K2").Value = "%Difference_Best"

# 2. Stale ZEUS numbers removed for p9 (row 11) and p11 (row 13)
$ws.Range("J11:K11").ClearContents()
$ws.Range("J13:K13").ClearContents()

# 3. pr9 (row 43) now has a ZEUS result
$ws.Range("J43").Value = 13729.7939453125
$ws.Range("K43").Value = 36.0

# 4. pr10 (row 44) ZEUS result removed
$ws.Range("J44:K44").ClearContents()

# 5. Columns G, H, J, K revert to default (non-custom) width
$ws.Columns.Item(7).ClearFormats()
$ws.Columns.Item(8).ClearFormats()
$ws.Columns.Item(10).ClearFormats()
$ws.Columns.Item(11).ClearFormats()

# 6. ClearFormats() on a whole column materialises blank placeholder
#    cells in every row that previously had none; sweep those back out
#    so we don't introduce cells the diff doesn't have.
for ($r = 1; $r -le 44; $r++) {
    foreach ($colIdx in 7,8,10,11) {
        $cell = $ws.Cells.Item($r, $colIdx)
        if ($cell.Text -eq "") {
            $cell.ClearContents()
        }
    }
}
